$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:E2").NumberFormat = "@"
$ws.Range("D3:E3").NumberFormat = "@"
$ws.Range("D4:E4").NumberFormat = "@"
$ws.Range("D5:E5").NumberFormat = "@"
$ws.Range("D6:E6").NumberFormat = "@"
$ws.Range("D7:E7").NumberFormat = "@"
$ws.Range("D8:E8").NumberFormat = "@"
$ws.Range("D9:E9").NumberFormat = "@"
$ws.Range("D10:E10").NumberFormat = "@"
$ws.Range("D11:E11").NumberFormat = "@"
$ws.Range("D12:E12").NumberFormat = "@"
$ws.Range("D13:E13").NumberFormat = "@"
$ws.Range("D14:E14").NumberFormat = "@"
$ws.Range("D15:E15").NumberFormat = "@"
$ws.Range("D16:E16").NumberFormat = "@"
$ws.Range("D17:E17").NumberFormat = "@"
$ws.Range("D18:E18").NumberFormat = "@"
$ws.Range("D19:E19").NumberFormat = "@"
$ws.Range("D20:E20").NumberFormat = "@"
$ws.Range("D21:E21").NumberFormat = "@"
$ws.Range("D22:E22").NumberFormat = "@"
$ws.Range("D23:E23").NumberFormat = "@"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("D25:E25").NumberFormat = "@"
$ws.Range("D26:E26").NumberFormat = "@"
$ws.Range("D27:E27").NumberFormat = "@"
$ws.Range("D39:E39").NumberFormat = "@"
$ws.Range("D40:E40").NumberFormat = "@"
$ws.Range("D41:E41").NumberFormat = "@"
$ws.Range("D42:E42").NumberFormat = "@"
$ws.Range("D43:E43").NumberFormat = "@"
$ws.Range("D44:E44").NumberFormat = "@"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("D46:E46").NumberFormat = "@"
$ws.Range("D47:E47").NumberFormat = "@"

$ws.Range('D2').Value = '304.65'
$ws.Range('E2').Value = '3.89%'
$ws.Range('D3').Value = '32.22'
$ws.Range('E3').Value = '4.82%'
$ws.Range('B4').Value = 'LEO'
$ws.Range('C4').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D4').Value = '3.443'
$ws.Range('E4').Value = '-1.13%'
$ws.Range('B5').Value = 'HuobiToken'
$ws.Range('C5').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D5').Value = '5.304'
$ws.Range('E5').Value = '2.95%'
$ws.Range('B6').Value = 'Cronos'
$ws.Range('C6').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D6').Value = '0.07586'
$ws.Range('E6').Value = '6.49%'
$ws.Range('B7').Value = 'KuCoinToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range('D7').Value = '7.868'
$ws.Range('E7').Value = '4.42%'
$ws.Range('B8').Value = 'GateToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D8').Value = '3.860'
$ws.Range('E8').Value = '6.29%'
$ws.Range('B9').Value = 'FTXToken'
$ws.Range('C9').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D9').Value = '1.648'
$ws.Range('E9').Value = '17.09%'
$ws.Range('B10').Value = 'MXToken'
$ws.Range('C10').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D10').Value = '0.9280'
$ws.Range('E10').Value = '1.21%'
$ws.Range('B11').Value = 'WazirX'
$ws.Range('C11').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D11').Value = '0.1692'
$ws.Range('E11').Value = '3.76%'
$ws.Range('B12').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C12').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D12').Value = '0.08017'
$ws.Range('E12').Value = '3.74%'
$ws.Range('B13').Value = 'MandalaExchangeToken'
$ws.Range('C13').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D13').Value = '0.08092'
$ws.Range('E13').Value = '4.30%'
$ws.Range('B14').Value = 'BitrueCoin'
$ws.Range('C14').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D14').Value = '0.03061'
$ws.Range('E14').Value = '3.92%'
$ws.Range('B15').Value = 'BitMartToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D15').Value = '0.09919'
$ws.Range('E15').Value = '10.17%'
$ws.Range('B16').Value = 'BitForexToken'
$ws.Range('C16').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D16').Value = '0.001492'
$ws.Range('E16').Value = '-6.38%'
$ws.Range('B17').Value = 'CoinExToken'
$ws.Range('C17').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('D17').Value = '0.04595'
$ws.Range('E17').Value = '1.47%'
$ws.Range('B18').Value = 'TigerCash'
$ws.Range('C18').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D18').Value = '0.006457'
$ws.Range('E18').Value = '0.02%'
$ws.Range('D19').Value = '2.234'
$ws.Range('E19').Value = '-0.44%'
$ws.Range('D20').Value = '0.3301'
$ws.Range('E20').Value = '1.55%'
$ws.Range('D21').Value = '0.1343'
$ws.Range('E21').Value = '-1.62%'
$ws.Range('D22').Value = '4.539'
$ws.Range('E22').Value = '17.98%'
$ws.Range('D23').Value = '0.1615'
$ws.Range('E23').Value = '0.91%'
$ws.Range('E24').Value = '0.21%'
$ws.Range('D25').Value = '0.004488'
$ws.Range('E25').Value = '5.95%'
$ws.Range('D26').Value = '0.0001394'
$ws.Range('E26').Value = '19.17%'
$ws.Range('D27').Value = '0.0001774'
$ws.Range('E27').Value = '4.96%'
$ws.Range('D39').Value = '0.01724'
$ws.Range('E39').Value = '2,522.73%'
$ws.Range('D40').Value = '0.04516'
$ws.Range('E40').Value = '2.38%'
$ws.Range('D41').Value = '0.006977'
$ws.Range('E41').Value = '-0.61%'
$ws.Range('D42').Value = '0.1360'
$ws.Range('E42').Value = '6.86%'
$ws.Range('B43').Value = 'CEJI'
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range('D43').Value = '0.002072'
$ws.Range('E43').Value = '-6.27%'
$ws.Range('B44').Value = 'LocalTraders'
$ws.Range('C44').Value = 'https://coinranking.com/coin/E6DwMU2zXb+localtraders-lct'
$ws.Range('D44').Value = '0.01378'
$ws.Range('E44').Value = '4.42%'
$ws.Range('E45').Value = '5.13%'
$ws.Range('D46').Value = '0.7191'
$ws.Range('E46').Value = '-62.73%'
$ws.Range('D47').Value = '0.01295'
$ws.Range('E47').Value = '-0.40%'

Write-Output "Applied all cell updates."
